$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 2; $r -le $rows; $r++) {
    $a = $ws.Cells.Item($r, 1)
    $b = $ws.Cells.Item($r, 2)

    if ($a.Value2 -eq "⬛") {
        $a.Value = "📘"
    }
    if ($b.Value2 -eq "noir") {
        $b.Value = "bleu"
    }
}
